# Add the 2022-Q4 sheet (new quarterly holdings data) and update the
# "总计" (totals) summary sheet to include it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new worksheet named "2022-Q4" right after "总计" so the
#    tab order becomes: 总计, 2022-Q4, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

# Header row (same headers used by the other quarterly sheets)
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"
$q4.Range("B1:H1").Font.Bold = $true
$q4.Range("B1:H1").HorizontalAlignment = -4108
$q4.Range("B1:H1").VerticalAlignment = -4160
$q4.Range("B1:H1").Borders.LineStyle = 1

# Data row: single fund holding for 2022-Q4
$q4.Range("A2").Value = 0
$q4.Range("A2").Font.Bold = $true
$q4.Range("A2").HorizontalAlignment = -4108
$q4.Range("A2").VerticalAlignment = -4160
$q4.Range("A2").Borders.LineStyle = 1

$q4.Range("B2").NumberFormat = "@"
$q4.Range("B2").Value = "007943"
$q4.Range("B2").NumberFormat = "General"

$q4.Range("C2").Value = "富安达中证 500 指数增强"

$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "0.30"
$q4.Range("D2").NumberFormat = "General"

$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "78.45"
$q4.Range("E2").NumberFormat = "General"

$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "1.11"
$q4.Range("F2").NumberFormat = "General"

$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "0.0033"
$q4.Range("G2").NumberFormat = "General"

$q4.Range("H2").Value = 9

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new data row for 2022-Q4 right
#    after the header row, pushing the existing quarters down by one
#    row and renumbering the index column (A). Values are written
#    directly (bottom row first) instead of using Rows.Insert so that
#    every pre-existing, already-styled cell is simply overwritten in
#    place and keeps its original formatting; only the brand new row
#    (row 6) needs its style applied explicitly.
# ---------------------------------------------------------------------

# New row 6 (was the 2021-Q4 row, now shifted down one row)
$total.Range("A6").Value = 4
$total.Range("A6").Font.Bold = $true
$total.Range("A6").HorizontalAlignment = -4108
$total.Range("A6").VerticalAlignment = -4160
$total.Range("A6").Borders.LineStyle = 1
$total.Range("B6").Value = "2021-Q4"
$total.Range("C6").Value = 3
$total.Range("D6").Value = 0.3

# Row 5 (was 2022-Q1)
$total.Range("A5").Value = 3
$total.Range("B5").Value = "2022-Q1"
$total.Range("C5").Value = 1
$total.Range("D5").Value = 0.17

# Row 4 (was 2022-Q2)
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.21

# Row 3 (was 2022-Q3)
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.13

# Row 2 (new 2022-Q4 row)
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0
